$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "release/1.0.1"
$ws.Range("B4").Value = "X"
$ws.Range("C4").Value = "X"
$ws.Range("D4").Value = "X"
